$d = $word.ActiveDocument

# Locate the "Full-Stack Development and Data Engineering" paragraph under the
# Siege Analytics / PARTNER role and insert four new bullet points right after it.
$rng = $d.Content
$found = $rng.Find.Execute("Full-Stack Development and Data Engineering", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)
    $newBullets = "`r• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States"
    $newBullets += "`r• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times"
    $newBullets += "`r• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis"
    $newBullets += "`r• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
    $rng.InsertAfter($newBullets)
}
